# edit.ps1 - applies "New crime data collected" changes to cs-en-us-007pct.xlsx
# Updates: volume/issue number, report date range, and weekly/28-day/YTD/2yr
# crime-complaint figures + percent changes for rows 14-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: issue number 28 -> 29, report week 7/7-7/13 -> 7/14-7/20 ---
$ws.Range("A8").Characters(21,2).Text = "29"
$ws.Range("C9").Characters(27,8).Text = "7/14/2025"
$ws.Range("C9").Characters(47,9).Text = "7/20/2025"

# --- Cells that become "N/A" (text) cells: copy style+text from an existing
#     identically-styled N/A cell so the shared-string/style stays canonical ---
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E31").PasteSpecial(-4163)

# --- Cells that change from "N/A" text to real numbers: set the number format
#     (reuses the existing numeric style) before assigning the value ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H15").Value = -100
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 4
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E20").Value = -75
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 1
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = -100
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 2

# --- Remaining cells: plain numeric value updates ---
$ws.Range("N14").Value = -71.428571428571
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -14.285714285714
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 70
$ws.Range("I16").Value = 71
$ws.Range("J16").Value = 57
$ws.Range("K16").Value = 24.561403508771
$ws.Range("L16").Value = -1.388888888888
$ws.Range("M16").Value = -8.974358974358
$ws.Range("N16").Value = -83.488372093023
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 121
$ws.Range("J17").Value = 129
$ws.Range("K17").Value = -6.201550387596
$ws.Range("L17").Value = 1.680672268907
$ws.Range("M17").Value = 45.783132530120
$ws.Range("N17").Value = -3.2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 50
$ws.Range("J18").Value = 79
$ws.Range("K18").Value = -36.708860759493
$ws.Range("L18").Value = -46.808510638297
$ws.Range("M18").Value = 16.279069767441
$ws.Range("N18").Value = -74.874371859296
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -22.916666666666
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 251
$ws.Range("K19").Value = -8.366533864541
$ws.Range("L19").Value = -29.230769230769
$ws.Range("M19").Value = 60.839160839160
$ws.Range("N19").Value = -8.366533864541
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -44.117647058823
$ws.Range("L20").Value = -36.666666666666
$ws.Range("M20").Value = -36.666666666666
$ws.Range("N20").Value = -90.821256038647
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -15.625
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = -6.976744186046
$ws.Range("I21").Value = 499
$ws.Range("J21").Value = 557
$ws.Range("K21").Value = -10.412926391382
$ws.Range("L21").Value = -23.112480739599
$ws.Range("M21").Value = 30.628272251308
$ws.Range("N21").Value = -59.430894308943
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = 19
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 72.727272727272
$ws.Range("L22").Value = 72.727272727272
$ws.Range("M22").Value = 111.111111111111
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 13.333333333333
$ws.Range("I23").Value = 94
$ws.Range("J23").Value = 119
$ws.Range("K23").Value = -21.008403361344
$ws.Range("L23").Value = 2.173913043478
$ws.Range("M23").Value = 11.904761904761
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = 4.347826086956
$ws.Range("I24").Value = 655
$ws.Range("J24").Value = 730
$ws.Range("K24").Value = -10.273972602739
$ws.Range("L24").Value = -3.534609720176
$ws.Range("M24").Value = 57.831325301204
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -13.043478260869
$ws.Range("F25").Value = 78
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = -2.5
$ws.Range("I25").Value = 376
$ws.Range("J25").Value = 496
$ws.Range("K25").Value = -24.193548387096
$ws.Range("L25").Value = 3.013698630136
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 150
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 227
$ws.Range("J26").Value = 233
$ws.Range("K26").Value = -2.575107296137
$ws.Range("L26").Value = -6.584362139917
$ws.Range("M26").Value = 20.105820105820
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -53.846153846153
$ws.Range("L27").Value = -53.846153846153
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("I28").Value = 30
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = 7.142857142857
$ws.Range("N29").Value = -83.333333333333
$ws.Range("N30").Value = -80
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = -50
$ws.Range("L31").Value = 66.666666666666
